# Add 2022-Q4 data
# ------------------------------------------------------------------
# Before: 总计 | 2022-Q2 | 2022-Q1
# After : 总计 | 2022-Q4 | 2022-Q2 | 2022-Q1
#
# A brand-new "2022-Q4" sheet (with its own fund holdings table) is
# inserted right after "总计" and before "2022-Q2". The "总计" summary
# sheet gets a new row for 2022-Q4 (pushing 2022-Q2 / 2022-Q1 down).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$sheetTotal = $wb.Worksheets.Item(1)   # 总计
$sheetQ2    = $wb.Worksheets.Item(2)   # 2022-Q2 (current position)
$sheetQ1    = $wb.Worksheets.Item(3)   # 2022-Q1 (current position)

# --------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by duplicating the "2022-Q2"
#    sheet (same column layout/styling) right after "总计", then
#    swap in the 2022-Q4 fund data and drop the extra 3rd data row.
# --------------------------------------------------------------
$sheetQ2.Copy([System.Type]::Missing, $sheetTotal)
$sheetQ4 = $wb.Worksheets.Item(2)
$sheetQ4.Name = "2022-Q4"

# The copied sheet has 3 data rows (like 2022-Q2); 2022-Q4 only needs 2.
$sheetQ4.Rows.Item(4).Delete()

# Columns B..G hold "numeric looking" text (fund codes / percentages)
# that must stay text (leading zeros / exact decimals) instead of
# being auto-converted to numbers.
$sheetQ4.Range("B2:G3").NumberFormat = "@"

$sheetQ4.Range("B2").Value = "005460"
$sheetQ4.Range("C2").Value = "银河嘉谊灵活配置混合C"
$sheetQ4.Range("D2").Value = "2.65"
$sheetQ4.Range("E2").Value = "39.19"
$sheetQ4.Range("F2").Value = "1.35"
$sheetQ4.Range("G2").Value = "0.0358"
$sheetQ4.Range("H2").Value = 1

$sheetQ4.Range("B3").Value = "005459"
$sheetQ4.Range("C3").Value = "银河嘉谊灵活配置混合A"
$sheetQ4.Range("D3").Value = "0.01"
$sheetQ4.Range("E3").Value = "39.19"
$sheetQ4.Range("F3").Value = "1.35"
$sheetQ4.Range("G3").Value = "0.0001"
$sheetQ4.Range("H3").Value = 1

# Drop the temporary text number-format again so the cells go back to
# the plain (unstyled) look used by the rest of the data rows.
$sheetQ4.Range("B2:G3").Style = "Normal"

# --------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert the 2022-Q4 row at
#    the top of the data and push 2022-Q2 / 2022-Q1 down by one row.
# --------------------------------------------------------------
$sheetTotal.Range("A3").Copy()
$sheetTotal.Range("A4").PasteSpecial(-4122)   # xlPasteFormats
$sheetTotal.Range("A4").Value = 2
$sheetTotal.Range("B4").Value = "2022-Q1"
$sheetTotal.Range("C4").Value = 3
$sheetTotal.Range("D4").Value = 0.05

$sheetTotal.Range("B3").Value = "2022-Q2"
$sheetTotal.Range("C3").Value = 3
$sheetTotal.Range("D3").Value = 0.05

$sheetTotal.Range("B2").Value = "2022-Q4"
$sheetTotal.Range("C2").Value = 2
$sheetTotal.Range("D2").Value = 0.04

# --------------------------------------------------------------
# 3. Restore the originally selected tab (2022-Q1, now the 4th sheet)
#    as the active / highlighted worksheet tab.
# --------------------------------------------------------------
$wb.Worksheets.Item(4).Activate()
